# edit.ps1 - applies the diff described changes to the document
$d = $word.ActiveDocument

# 1) Merge the runs describing "но не смог зарегистрировать..." (removes proofErr wrapped spell/gram-check runs)
$d.Content.Find.Execute(
    "но не смог зарегистрировать аккаунт. Написал тикет в поддержку, они ответили и вместе проблему мы решили. Немножко поизучал гайд об начале работы с этим сервисом.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "но не смог зарегистрировать аккаунт. Написал тикет в поддержку, они ответили и вместе проблему мы решили. Немножко поизучал гайд об начале работы с этим сервисом.",
    2) | Out-Null

# 2) Merge the runs describing "которое будет использовать распознавание лиц и будет иметь возможность распознавать клиентов."
$d.Content.Find.Execute(
    "которое будет использовать распознавание лиц и будет иметь возможность распознавать клиентов.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "которое будет использовать распознавание лиц и будет иметь возможность распознавать клиентов.",
    2) | Out-Null

# 3) Merge the runs describing "закончил тестовый вариант..." (keep the leading space)
$d.Content.Find.Execute(
    " закончил тестовый вариант, все работает, но пока что не придумал как его поменять под то, что необходимо.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " закончил тестовый вариант, все работает, но пока что не придумал как его поменять под то, что необходимо.",
    2) | Out-Null

# 4) Merge the runs describing "Продолжил изучать возможности локальной проверки...Лауриса."
$d.Content.Find.Execute(
    "Продолжил изучать возможности локальной проверки на наличие лиц на картинке, ввел в курс дела практиканта Лауриса.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Продолжил изучать возможности локальной проверки на наличие лиц на картинке, ввел в курс дела практиканта Лауриса.",
    2) | Out-Null

# 5) Restructure the tail: split "21.09.2020"/"22.09.2020" paragraphs, add a new "16.09.2020" entry
#    before 21.09.2020, add a new "22.09.2020" entry (dlib paragraph) and make the bookmark-only
#    paragraph the new last (empty) paragraph.
$p21 = $null
$p22 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs.Item($i).Range.Text
    if ($ptxt.StartsWith("21.09.2020")) { $p21 = $d.Paragraphs.Item($i) }
    if ($ptxt.StartsWith("22.09.2020")) { $p22 = $d.Paragraphs.Item($i) }
}

$tailRange = $d.Range($p21.Range.Start, $p22.Range.End)
$tailRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:lang w:val="lv-LV"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">16.09.2020: </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Вместе с Лаурисом составили примерный план работы и к ней приступили, я начал заниматься внедрением программы для разделения видео на кадры в </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>C</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="lv-LV"/></w:rPr><w:t>#.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>21.09.2020:</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:lang w:val="lv-LV"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Работал над частью, которая должна локально проверить, есть ли на конкретной фотографии лицо человека.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>22.09.2020:</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Пытался найти возможности для улучшения качества проверки на наличие лиц в библиотеке </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>dlib</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
